$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) is treated as text so values like "1.000",
# "0.000007932", trailing-zero decimals, etc. are preserved verbatim
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.231.77'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.861.42'
$ws.Range("E3").Value = '  +0.69%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7096'
$ws.Range("E5").Value = '  +0.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '237.86'
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.08177'
$ws.Range("E8").Value = '  +10.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3042'
$ws.Range("E9").Value = '  -0.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.27'
$ws.Range("E10").Value = '  -0.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08166'
$ws.Range("E11").Value = '  +0.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.877.61'
$ws.Range("E12").Value = '  +2.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.172'
$ws.Range("E13").Value = '  -0.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7086'
$ws.Range("E14").Value = '  -2.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.73'
$ws.Range("E15").Value = '  +1.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.236.10'
$ws.Range("E16").Value = '  +0.24%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007932'
$ws.Range("E17").Value = '  +3.80%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.785'
$ws.Range("E18").Value = '  +0.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.36'
$ws.Range("E19").Value = '  +2.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.13'
$ws.Range("E20").Value = '  -0.60%  '
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.102.81'
$ws.Range("E22").Value = '  +0.40%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.404'
$ws.Range("E24").Value = '  -2.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.73'
$ws.Range("E25").Value = '  +1.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.963'
$ws.Range("E26").Value = '  -0.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1453'
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.08'
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.960'
$ws.Range("E29").Value = '  -0.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.425'
$ws.Range("E30").Value = '  +1.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.487'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.392'
$ws.Range("E32").Value = '  -3.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.028'
$ws.Range("E33").Value = '  +0.96%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05228'
$ws.Range("E34").Value = '  +0.65%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.170'
$ws.Range("E35").Value = '  -1.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7072'
$ws.Range("E36").Value = '  +0.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.001'
$ws.Range("E37").Value = '  -3.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.674'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01856'
$ws.Range("E39").Value = '  -0.77%  '
$ws.Range("E40").Value = '  +2.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.145.60'
$ws.Range("E41").Value = '  +6.82%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9238'
$ws.Range("E42").Value = '  -3.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4285'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.890'
$ws.Range("E44").Value = '  -2.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.38'
$ws.Range("E45").Value = '  +0.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.83'
$ws.Range("E47").Value = '  +0.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.775'
$ws.Range("E48").Value = '  +1.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.001.12'
$ws.Range("E49").Value = '  +0.67%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.215'
$ws.Range("E50").Value = '  +1.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.956'
$ws.Range("E51").Value = '  -1.43%  '
